$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a numeric-looking value while preserving it as text (matches
# original inlineStr/text cell content) and without leaving stray styling
# (NumberFormat/quotePrefix) behind on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Price (column D) updates ---
Set-TextValue "D2" "247.34"
Set-TextValue "D3" "22.40"
Set-TextValue "D4" "5.471"
Set-TextValue "D5" "0.05623"
Set-TextValue "D6" "6.459"
Set-TextValue "D7" "0.8044"
Set-TextValue "D9" "0.1419"
Set-TextValue "D10" "0.07301"
Set-TextValue "D12" "0.02922"
Set-TextValue "D13" "0.09245"
Set-TextValue "D14" "0.001668"
Set-TextValue "D15" "3.212"
Set-TextValue "D16" "0.04748"
Set-TextValue "D17" "0.0005817"
$ws.Range("E17").Value = "16OneONEWorstin24h"
Set-TextValue "D18" "0.006364"
Set-TextValue "D20" "0.001058"
Set-TextValue "D22" "3.983"
Set-TextValue "D23" "3.378"
Set-TextValue "D24" "2.122"
Set-TextValue "D27" "0.0003304"
Set-TextValue "D40" "0.04153"

# --- Row 41 / Row 43 swap (BKEXToken <-> KickToken) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006901"
$ws.Range("E41").Value = "40KickTokenKICK"

Set-TextValue "D42" "0.003505"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1038"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.009152"
Set-TextValue "D45" "0.00005659"
Set-TextValue "D47" "0.6808"
Set-TextValue "D48" "0.01591"
$ws.Range("E48").Value = "47BOLOBOLO"
Set-TextValue "D49" "0.00002103"
Set-TextValue "D50" "0.01011"
